# Waste model update: revise the "Methane recovery in landfill to generate
# electricity" demand share ramp (row 3) and tidy up the number formatting
# of the "Total solid waste" row (row 2) to match the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand_Projection")

# --- Row 2 (E5TSWTSW / Total solid waste): apply the 4-decimal number
#     format used everywhere else on the sheet to the I2:AG2 run (AH2
#     already had it).
$ws.Range("I2:AG2").NumberFormat = "0.0000"

# --- Row 3 (E5TSWLANDFILL_ELEC / Methane recovery in landfill): the ramp
#     now climbs twice as fast (0.006 step instead of 0.003) and, after
#     peaking at 0.03 in 2025 (column X), drops back to 0 instead of
#     plateauing at 0.03 through 2035.
$ws.Range("T3").Value = 0.006
$ws.Range("U3").Value = 0.012
$ws.Range("V3").Value = 0.018
$ws.Range("W3").Value = 0.024
$ws.Range("X3").Value = 0.03
$ws.Range("Y3:AH3").Value = 0

# --- Reflect the user's final selection (they had just edited/selected
#     the row-3 figures) instead of the original I2 single-cell selection.
$ws.Range("I3:AH3").Select()
